$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.613.12"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.582.45"
$ws.Range("E3").Value = "  -3.11%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.42"
$ws.Range("E5").Value = "  -2.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.502"
$ws.Range("E6").Value = "  -3.30%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.21"
$ws.Range("E8").Value = "  -5.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.253"
$ws.Range("E9").Value = "  -1.79%  "
$ws.Range("E10").Value = "  -3.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0866"
$ws.Range("E11").Value = "  -1.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.806.77"
$ws.Range("E12").Value = "  -3.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.594.89"
$ws.Range("E13").Value = "  -2.40%  "
$ws.Range("E14").Value = "  -4.43%  "
$ws.Range("E15").Value = "  -6.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.557.60"
$ws.Range("E16").Value = "  -1.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.77"
$ws.Range("E17").Value = "  -4.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.43"
$ws.Range("E18").Value = "  -5.07%  "
$ws.Range("E19").Value = "  -4.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0$([char]0x2083)0692"
$ws.Range("E20").Value = "  -3.71%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("E22").Value = "  -4.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.50"
$ws.Range("E23").Value = "  -5.71%  "
$ws.Range("E24").Value = "  -4.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.54"
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.70"
$ws.Range("E27").Value = "  -2.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.05"
$ws.Range("E28").Value = "  -3.17%  "
$ws.Range("E29").Value = "  -4.57%  "
$ws.Range("E30").Value = "  -2.56%  "
$ws.Range("E31").Value = "  -3.81%  "
$ws.Range("E32").Value = "  -5.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.371.50"
$ws.Range("E33").Value = "  -1.61%  "
$ws.Range("E34").Value = "  -5.35%  "
$ws.Range("E35").Value = "  -5.11%  "
$ws.Range("E36").Value = "  -6.46%  "
$ws.Range("E37").Value = "  -2.21%  "
$ws.Range("E38").Value = "  -4.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.535"
$ws.Range("E39").Value = "  -4.34%  "
$ws.Range("E40").Value = "  -4.27%  "
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("E42").Value = "  -3.75%  "
$ws.Range("E43").Value = "  -2.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.18"
$ws.Range("E44").Value = "  +1.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.56"
$ws.Range("E45").Value = "  -3.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.26"
$ws.Range("E46").Value = "  -3.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.716.23"
$ws.Range("E47").Value = "  -3.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.42"
$ws.Range("E48").Value = "  -1.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0$([char]0x2087)0999"
$ws.Range("E49").Value = "  -3.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0972"
$ws.Range("E50").Value = "  -4.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0497"
$ws.Range("E51").Value = "  -1.41%  "
